# "different logic for buyer's market choice"
# Update the Configuration sheet's simulation parameters and add a new
# MEMORY parameter row, then adjust the markets sheet view (freeze the
# first column / scroll) to match the author's new working state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Configuration sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Configuration")
$ws1.Select()

# PERIODS: 30 -> 50
$ws1.Range("B1").Value = 50

# AGENTS: 100 -> 4
$ws1.Range("B2").Value = 4

# CONTACTS: 17 -> 2, and pick up the "Calibri Light" style used by A3
# (copy formats only so no new style slot is created)
$ws1.Range("A3").Copy()
$ws1.Range("B3").PasteSpecial(-4122)
$ws1.Range("B3").Value = 2

# GUI: 0 -> 1
$ws1.Range("B9").Value = 1

# New row 11: MEMORY = -1, styled like the other label cells (A10)
$ws1.Range("A10").Copy()
$ws1.Range("A11").PasteSpecial(-4122)
$ws1.Range("A11").Value = "MEMORY"
$ws1.Range("B11").Value = -1

$ws1.Range("B2").Select()

# ---------------------------------------------------------------------
# markets sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("markets")
$ws2.Select()

# Freeze the first column, then move the selection back to A2.
$ws2.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("A2").Select()

# Restore Configuration as the active sheet/tab.
$ws1.Select()
$ws1.Range("B2").Select()
